# Add a "metadata" worksheet after the existing "data" worksheet and
# populate it with panel metadata, matching the target revision of the
# workbook (commit: "Refined metadata to be additional tab").
#
# Note: the F column ("time_taken") timestamps on the "data" sheet were
# also refreshed to a later query run; update those too.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Update the "time_taken" timestamps on the data sheet (column F,
#    rows 2-19) to reflect the newer panel query.
# ---------------------------------------------------------------------
$timeTaken = @(
    "2021-10-05 14:21:49.760495",
    "2021-10-05 14:21:49.760503",
    "2021-10-05 14:21:49.760506",
    "2021-10-05 14:21:49.760509",
    "2021-10-05 14:21:49.760512",
    "2021-10-05 14:21:49.760514",
    "2021-10-05 14:21:49.760517",
    "2021-10-05 14:21:49.760520",
    "2021-10-05 14:21:49.760522",
    "2021-10-05 14:21:49.760525",
    "2021-10-05 14:21:49.760528",
    "2021-10-05 14:21:49.760531",
    "2021-10-05 14:21:49.760533",
    "2021-10-05 14:21:49.760536",
    "2021-10-05 14:21:49.760538",
    "2021-10-05 14:21:49.760541",
    "2021-10-05 14:21:49.760544",
    "2021-10-05 14:21:49.760546"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the bold/bordered header style used on the data sheet (style
# applied to data!B1) onto the metadata header row and the index cell.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row.
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row.
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Neurotransmitter disorders"
$metaSheet.Cells.Item(2, 3).Value = 219
$metaSheet.Cells.Item(2, 4).Value = "'1.9"
$metaSheet.Cells.Item(2, 5).Value = "2021-04-01T15:08:40.081474Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:21:49.756830"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/219/?format=json"

$dataSheet.Activate()
